# Change the table style on the three data tables (slides 14-16) from the
# presentation's custom default table style to the built-in PowerPoint
# table style "{304FF40B-6C6E-498B-AD1D-2C46EE2DC3AB}" (equivalent to
# selecting a different style in the Table Design > Table Styles gallery).

$p = $ppt.ActivePresentation
$newStyleId = "{304FF40B-6C6E-498B-AD1D-2C46EE2DC3AB}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
